$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.642.47'
$ws.Range('D2').Style = $ws.Range('B2').Style
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.805.71'
$ws.Range('D3').Style = $ws.Range('B2').Style
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.76'
$ws.Range('D5').Style = $ws.Range('B2').Style
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5462'
$ws.Range('D7').Style = $ws.Range('B2').Style
$ws.Range('E7').Value = '  -4.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3808'
$ws.Range('D8').Style = $ws.Range('B2').Style
$ws.Range('E8').Value = '  -1.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07517'
$ws.Range('D9').Style = $ws.Range('B2').Style
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.40'
$ws.Range('D10').Style = $ws.Range('B2').Style
$ws.Range('E10').Value = '  -1.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.115'
$ws.Range('D11').Style = $ws.Range('B2').Style
$ws.Range('E11').Value = '  -2.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = $ws.Range('B2').Style
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.73'
$ws.Range('D13').Style = $ws.Range('B2').Style
$ws.Range('E13').Value = '  -2.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.166'
$ws.Range('D14').Style = $ws.Range('B2').Style
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.400'
$ws.Range('D15').Style = $ws.Range('B2').Style
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.790.19'
$ws.Range('D16').Style = $ws.Range('B2').Style
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '90.29'
$ws.Range('D17').Style = $ws.Range('B2').Style
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001068'
$ws.Range('D18').Style = $ws.Range('B2').Style
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06483'
$ws.Range('D19').Style = $ws.Range('B2').Style
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.34'
$ws.Range('D21').Style = $ws.Range('B2').Style
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.944'
$ws.Range('D22').Style = $ws.Range('B2').Style
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.622.33'
$ws.Range('D23').Style = $ws.Range('B2').Style
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.12'
$ws.Range('D24').Style = $ws.Range('B2').Style
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.103'
$ws.Range('D25').Style = $ws.Range('B2').Style
$ws.Range('E25').Value = '  -1.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.68'
$ws.Range('D26').Style = $ws.Range('B2').Style
$ws.Range('E26').Value = '  +1.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.46'
$ws.Range('D27').Style = $ws.Range('B2').Style
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.365'
$ws.Range('D28').Style = $ws.Range('B2').Style
$ws.Range('E28').Value = '  -3.62%  '
$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.002.38'
$ws.Range('D29').Style = $ws.Range('B2').Style
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.25'
$ws.Range('D30').Style = $ws.Range('B2').Style
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.117'
$ws.Range('D31').Style = $ws.Range('B2').Style
$ws.Range('E31').Value = '  -3.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1057'
$ws.Range('D32').Style = $ws.Range('B2').Style
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.644'
$ws.Range('D33').Style = $ws.Range('B2').Style
$ws.Range('E33').Value = '  -2.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.684'
$ws.Range('D34').Style = $ws.Range('B2').Style
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.06658'
$ws.Range('D35').Style = $ws.Range('B2').Style
$ws.Range('E35').Value = '  +8.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2258'
$ws.Range('D36').Style = $ws.Range('B2').Style
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02303'
$ws.Range('D37').Style = $ws.Range('B2').Style
$ws.Range('E37').Value = '  -0.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.777'
$ws.Range('D38').Style = $ws.Range('B2').Style
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.033'
$ws.Range('D39').Style = $ws.Range('B2').Style
$ws.Range('E39').Value = '  -0.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6252'
$ws.Range('D40').Style = $ws.Range('B2').Style
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.28'
$ws.Range('D41').Style = $ws.Range('B2').Style
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.197'
$ws.Range('D42').Style = $ws.Range('B2').Style
$ws.Range('E42').Value = '  +2.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.438'
$ws.Range('D43').Style = $ws.Range('B2').Style
$ws.Range('E43').Value = '  +4.33%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.25'
$ws.Range('D44').Style = $ws.Range('B2').Style
$ws.Range('E44').Value = '  -1.57%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5861'
$ws.Range('D45').Style = $ws.Range('B2').Style
$ws.Range('E45').Value = '  -2.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.696'
$ws.Range('D46').Style = $ws.Range('B2').Style
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '126.80'
$ws.Range('D47').Style = $ws.Range('B2').Style
$ws.Range('E47').Value = '  +3.67%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.950'
$ws.Range('D48').Style = $ws.Range('B2').Style
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.160'
$ws.Range('D49').Style = $ws.Range('B2').Style
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06890'
$ws.Range('D50').Style = $ws.Range('B2').Style
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.33'
$ws.Range('D51').Style = $ws.Range('B2').Style
$ws.Range('E51').Value = '  -1.19%  '
